$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Fecha" (A4) and "Fecha Pago" (L4) text values for row 4
# (L4 first so the shared-string table slot order matches the target file)
$ws.Range("L4").Value = ".12.15.2022"
$ws.Range("A4").Value = ".13.12.2022"

# Update the "Monto" column values for rows 2-4
$ws.Range("I2").Value = 500.5
$ws.Range("I3").Value = 5.5
$ws.Range("I4").Value = 1500.9

# The wider "Monto" values no longer fit the old best-fit column width,
# so widen column I (Monto) to fit the new content
$ws.Columns.Item(9).ColumnWidth = 7

# Update the active cell selection
$ws.Range("I5").Select()
